{"js": "const pairs = [\n  [\"2026-01-22 Thursday\", \"2026-01-23 Friday\"],\n  [\"57\u00d785=\", \"62\u00d795=\"],\n  [\"87\u00d764=\", \"26\u00d718=\"],\n  [\"59\u00d787=\", \"36\u00d754=\"],\n  [\"47\u00d736=\", \"93\u00d799=\"],\n  [\"14\u00d788=\", \"54\u00d719=\"],\n  [\"75\u00d799=\", \"31\u00d746=\"],\n  [\"52\u00d750=\", \"83\u00d779=\"],\n  [\"89\u00d786=\", \"77\u00d745=\"],\n  [\"49\u00d765=\", \"15\u00d717=\"],\n  [\"41\u00d785=\", \"59\u00d768=\"],\n  [\"65\u00d760=\", \"23\u00d740=\"],\n  [\"77\u00d796=\", \"98\u00d794=\"],\n  [\"63\u00d777=\", \"91\u00d773=\"],\n  [\"24\u00d751=\", \"50\u00d717=\"],\n  [\"99\u00d795=\", \"92\u00d739=\"],\n  [\"36\u00d750=\", \"65\u00d731=\"],\n  [\"40\u00d724=\", \"23\u00d780=\"],\n  [\"45\u00d725=\", \"26\u00d782=\"],\n  [\"73\u00d784=\", \"89\u00d785=\"],\n  [\"18\u00d734=\", \"52\u00d734=\"],\n  [\"95\u00d766=\", \"44\u00d788=\"],\n  [\"62\u00d756=\", \"68\u00d779=\"],\n  [\"70\u00d765=\", \"71\u00d738=\"],\n  [\"84\u00d732=\", \"62\u00d769=\"],\n  [\"97\u00d759=\", \"98\u00d786=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@('2026-01-22 Thursday', '2026-01-23 Friday')\n    ,@('57\u00d785=', '62\u00d795=')\n    ,@('87\u00d764=', '26\u00d718=')\n    ,@('59\u00d787=', '36\u00d754=')\n    ,@('47\u00d736=', '93\u00d799=')\n    ,@('14\u00d788=', '54\u00d719=')\n    ,@('75\u00d799=', '31\u00d746=')\n    ,@('52\u00d750=', '83\u00d779=')\n    ,@('89\u00d786=', '77\u00d745=')\n    ,@('49\u00d765=', '15\u00d717=')\n    ,@('41\u00d785=', '59\u00d768=')\n    ,@('65\u00d760=', '23\u00d740=')\n    ,@('77\u00d796=', '98\u00d794=')\n    ,@('63\u00d777=', '91\u00d773=')\n    ,@('24\u00d751=', '50\u00d717=')\n    ,@('99\u00d795=', '92\u00d739=')\n    ,@('36\u00d750=', '65\u00d731=')\n    ,@('40\u00d724=', '23\u00d780=')\n    ,@('45\u00d725=', '26\u00d782=')\n    ,@('73\u00d784=', '89\u00d785=')\n    ,@('18\u00d734=', '52\u00d734=')\n    ,@('95\u00d766=', '44\u00d788=')\n    ,@('62\u00d756=', '68\u00d779=')\n    ,@('70\u00d765=', '71\u00d738=')\n    ,@('84\u00d732=', '62\u00d769=')\n    ,@('97\u00d759=', '98\u00d786=')\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2)\n}"}
